# The deck ships with two theme parts:
#   ppt/theme/theme1.xml  -> bound to the slide master (currently the
#                             custom "Integral" colour scheme)
#   ppt/theme/theme2.xml  -> bound to the notes master (the default
#                             "Office" colour scheme: 44546A/E7E6E6/
#                             5B9BD5/ED7D31/A5A5A5/FFC000/4472C4/70AD47/
#                             0563C1/954F72)
#
# The target edit swaps which colours are applied to the slide master,
# i.e. the slide master ends up using the classic "Office" 12-colour
# theme palette (font/format schemes are already identical between the
# two themes, so only the colour scheme needs to change).
#
# PowerPoint's object model doesn't expose "replace this theme part
# wholesale"; the supported automation surface for recolouring a theme
# is ThemeColorScheme.Colors(i).RGB, so set every slot individually.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# Index order: 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10 accent1-6, 11 hlink, 12 folHlink
$tcs.Colors(1).RGB  = 0         # dk1      000000
$tcs.Colors(2).RGB  = 16777215  # lt1      FFFFFF
$tcs.Colors(3).RGB  = 6968388   # dk2      44546A
$tcs.Colors(4).RGB  = 15132391  # lt2      E7E6E6
$tcs.Colors(5).RGB  = 13998939  # accent1  5B9BD5
$tcs.Colors(6).RGB  = 3243501   # accent2  ED7D31
$tcs.Colors(7).RGB  = 10855845  # accent3  A5A5A5
$tcs.Colors(8).RGB  = 49407     # accent4  FFC000
$tcs.Colors(9).RGB  = 12874308  # accent5  4472C4
$tcs.Colors(10).RGB = 4697456   # accent6  70AD47
$tcs.Colors(11).RGB = 12673797  # hlink    0563C1
$tcs.Colors(12).RGB = 7491477   # folHlink 954F72
